# NIT-9002540882.xlsx - "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the totals at the top of the statement.
$ws.Range("E11").Value = 7773595
$ws.Range("F13").Value = 31

# 2) Insert a new detail row (period 2509) right after the current last
#    "2508" row for MARTHA GARCIA ABONDANO (row 46), pushing everything
#    below it down by one row. Copy formatting down from the row above
#    so the new row matches the existing table styling.
$ws.Rows("47:47").Insert()
$ws.Range("B46:J46").Copy()
$ws.Range("B47:J47").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B47").Value = "CC"
$ws.Range("C47").Value = "1020759457"
$ws.Range("D47").Value = "MARTHA GARCIA ABONDANO"
$ws.Range("E47").Value = "2509"
$ws.Range("F47").Value = 240000
$ws.Range("G47").Value = 6000000

# 3) The two rows that used to be the last two entries (FABIO and
#    MAURICIO, period 2508) are now one row further down (48 and 49) and
#    their period moves on to 2509 as well.
$ws.Range("E48").Value = "2509"
$ws.Range("E49").Value = "2509"

# 4) The row inserted in step 2 already pushed the signature block down
#    by one row (from 53/54 to 54/55), matching the target layout.

Write-Output "done"
